# Swap the presentation's design from the custom "Integral" theme back to
# the default "Office Theme" palette.
#
# The underlying theme part (ppt/theme/theme2.xml, the theme referenced by
# the presentation / slide master) is reachable through the
# ThemeColorScheme object exposed on the slide master's Theme. Re-pointing
# every theme colour slot (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink) to the
# stock Office palette reproduces the colours the default "Office Theme"
# uses (the font scheme and format scheme were already identical between
# the two themes in this deck, so nothing else needs to change).
#
# RGB() isn't available in this host, so the values below are the usual
# VBA-style BGR-packed long (R + G*256 + B*65536) for each hex colour:
#   1  dk1       000000 -> 0
#   2  lt1       FFFFFF -> 16777215
#   3  dk2       44546A -> 6968388
#   4  lt2       E7E6E6 -> 15132391
#   5  accent1   5B9BD5 -> 13998939
#   6  accent2   ED7D31 -> 3243501
#   7  accent3   A5A5A5 -> 10855845
#   8  accent4   FFC000 -> 49407
#   9  accent5   4472C4 -> 12874308
#   10 accent6   70AD47 -> 4697456
#   11 hlink     0563C1 -> 12673797
#   12 folHlink  954F72 -> 7491477

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$tcs = $m.Theme.ThemeColorScheme

$tcs.Colors(1).RGB  = 0
$tcs.Colors(2).RGB  = 16777215
$tcs.Colors(3).RGB  = 6968388
$tcs.Colors(4).RGB  = 15132391
$tcs.Colors(5).RGB  = 13998939
$tcs.Colors(6).RGB  = 3243501
$tcs.Colors(7).RGB  = 10855845
$tcs.Colors(8).RGB  = 49407
$tcs.Colors(9).RGB  = 12874308
$tcs.Colors(10).RGB = 4697456
$tcs.Colors(11).RGB = 12673797
$tcs.Colors(12).RGB = 7491477
